# Updates the cryptos price table (Price/Volume(1h) columns, plus two
# ranking swaps) to the latest scraped snapshot.
#
# Note: several "Price" values are digit strings that Excel's COM layer
# would otherwise auto-coerce into Doubles (losing trailing zeros / adding
# floating point noise, e.g. "602.10" -> 602.10000000000002). For those we
# assign with a leading apostrophe (forces text entry, like typing '602.10
# into the cell) and then reset Style to 'Normal' so the cell doesn't keep
# a lingering "quote prefix" number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.182.88'
$ws.Range('E2').Value = '  +0.76%  '

$ws.Range('D3').Value = '3.755.09'
$ws.Range('E3').Value = '  +0.94%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = "'602.10"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.13%  '

$ws.Range('D6').Value = "'167.17"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.17%  '

$ws.Range('D7').Value = '3.752.11'
$ws.Range('E7').Value = '  +0.83%  '

$ws.Range('E8').Value = '  +0.06%  '

$ws.Range('E9').Value = '  +1.07%  '

$ws.Range('E10').Value = '  +5.13%  '

$ws.Range('D11').Value = "'6.40"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.13%  '

$ws.Range('D12').Value = "'0.462"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.16%  '

$ws.Range('D13').Value = "'38.17"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.39%  '

$ws.Range('D14').Value = "'0.0000249"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.00%  '

$ws.Range('D15').Value = '4.385.24'
$ws.Range('E15').Value = '  +1.00%  '

$ws.Range('D16').Value = '3.736.31'
$ws.Range('E16').Value = '  +0.49%  '

$ws.Range('D17').Value = '69.282.68'
$ws.Range('E17').Value = '  +0.97%  '

$ws.Range('D18').Value = "'7.52"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.28%  '

$ws.Range('D19').Value = "'17.45"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.12%  '

$ws.Range('E20').Value = '  -1.54%  '

$ws.Range('E21').Value = '  +7.55%  '

$ws.Range('D22').Value = "'493.14"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.70%  '

$ws.Range('D23').Value = "'0.732"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.01%  '

$ws.Range('D24').Value = "'0.0000152"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.84%  '

$ws.Range('D25').Value = "'85.15"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.12%  '

$ws.Range('E26').Value = '  -0.59%  '

$ws.Range('D27').Value = "'12.46"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.32%  '

$ws.Range('E28').Value = '  -0.62%  '

$ws.Range('E29').Value = '  -0.03%  '

$ws.Range('D30').Value = "'8.32"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.65%  '

$ws.Range('E31').Value = '  +0.78%  '

$ws.Range('E32').Value = '  -3.81%  '

$ws.Range('D33').Value = "'31.70"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.79%  '

$ws.Range('D34').Value = '3.902.22'
$ws.Range('E34').Value = '  +0.96%  '

$ws.Range('B35').Value = 'RenzoRestakedETH'
$ws.Range('C35').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D35').Value = '3.687.48'
$ws.Range('E35').Value = '  +0.95%  '

$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = "'0.109"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.38%  '

$ws.Range('E37').Value = '  +3.10%  '

$ws.Range('D38').Value = "'0.141"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.25%  '

$ws.Range('E39').Value = '  +0.77%  '

$ws.Range('D40').Value = "'3.18"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.32%  '

$ws.Range('D41').Value = "'1.00"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.05%  '

$ws.Range('E42').Value = '  +0.67%  '

$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = "'2.02"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.18%  '

$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = "'48.64"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.17%  '

$ws.Range('D45').Value = "'426.82"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.44%  '

$ws.Range('D46').Value = "'8.50"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.00%  '

$ws.Range('E47').Value = '  +0.00%  '

$ws.Range('D48').Value = "'40.30"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.62%  '

$ws.Range('D49').Value = "'141.37"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.02%  '

$ws.Range('D50').Value = '2.798.44'
$ws.Range('E50').Value = '  +1.40%  '

$ws.Range('D51').Value = "'0.0354"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.19%  '
